# Update data: 18 July 2022
# Adds the newest month (1 June 2022, serial 44713) of unemployment data
# to the "Canada" sheet and the "Province" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Canada" (sheet1): append row 31
# ---------------------------------------------------------------------
$wsCanada = $wb.Worksheets.Item("Canada")

$wsCanada.Range("A31").Value = 44713
$wsCanada.Range("A31").NumberFormat = "d-mmm-yy"
$wsCanada.Range("B31").Value = "Canada"
$wsCanada.Range("B31").NumberFormat = "d-mmm-yy"
$wsCanada.Range("D31").Value = 1003.5
$wsCanada.Range("E31").Value = 1127
$wsCanada.Range("C31").Formula = "=(D31-E31)/E31*100"

$wsCanada.Range("A31").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "Province" (sheet2): append rows 292-301 (one per province)
# ---------------------------------------------------------------------
$wsProvince = $wb.Worksheets.Item("Province")

$provinceRows = @(
    @{ Row = 292; Name = "Newfoundland & Labrador"; D = 25.2;   E = 33.9;  DateStyle = $true  },
    @{ Row = 293; Name = "Prince Edward Island";     D = 4.5;   E = 7.5;   DateStyle = $false },
    @{ Row = 294; Name = "Nova Scotia";              D = 35.9;  E = 34.5;  DateStyle = $false },
    @{ Row = 295; Name = "New Brunswick";            D = 24;    E = 32.3;  DateStyle = $false },
    @{ Row = 296; Name = "Quebec";                   D = 193.9; E = 223.7; DateStyle = $false },
    @{ Row = 297; Name = "Ontario";                  D = 415.1; E = 433.1; DateStyle = $false },
    @{ Row = 298; Name = "Manitoba";                 D = 26.5;  E = 39;    DateStyle = $false },
    @{ Row = 299; Name = "Saskatchewan";              D = 23.6;  E = 31.5;  DateStyle = $false },
    @{ Row = 300; Name = "Alberta";                  D = 122.5; E = 165.3; DateStyle = $false },
    @{ Row = 301; Name = "British Columbia";         D = 132.2; E = 126.2; DateStyle = $false }
)

foreach ($r in $provinceRows) {
    $row = $r.Row

    $wsProvince.Range("A$row").Value = 44713
    $wsProvince.Range("A$row").NumberFormat = "d-mmm-yy"

    $wsProvince.Range("B$row").Value = $r.Name
    if ($r.DateStyle) {
        $wsProvince.Range("B$row").NumberFormat = "d-mmm-yy"
    }

    $wsProvince.Range("D$row").Value = $r.D
    $wsProvince.Range("E$row").Value = $r.E
    $wsProvince.Range("C$row").Formula = "=(D$row-E$row)/E$row*100"
}

$wsProvince.Range("D302").Select() | Out-Null
